$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$NL = [char]10

# --- Skill Level ("D") column updates: Trained/Aware -> Proficient ---
$ws.Range("D4").Value  = "Proficient"
$ws.Range("D7").Value  = "Proficient"
$ws.Range("D9").Value  = "Proficient"
$ws.Range("D11").Value = "Proficient"
$ws.Range("D14").Value = "Proficient"
$ws.Range("D20").Value = "Proficient"
$ws.Range("D22").Value = "Proficient"
$ws.Range("D23").Value = "Proficient"
$ws.Range("D25").Value = "Proficient"
$ws.Range("D32").Value = "Proficient"

# --- Evidence ("E") column updates ---
$ws.Range("E2").Value = "> Pearson Language test" + $NL + "> Liaison with various parties" + $NL + "> PoE"
$ws.Range("E4").Value = "> PoE " + $NL + "> OneNote Notebooks"
$ws.Range("E14").Value = "> Worked in critical services environment"

$elevenYearsText = "> 11 years in team envirnoment" + $NL + "> PoE"
$ws.Range("E15").Value = $elevenYearsText
$ws.Range("E16").Value = $elevenYearsText
$ws.Range("E17").Value = $elevenYearsText
$ws.Range("E18").Value = $elevenYearsText

# --- Row heights (grow to fit the new multi-line evidence text) ---
$ws.Rows.Item(2).RowHeight = 45
$ws.Rows.Item(15).RowHeight = 45.75
$ws.Rows.Item(16).RowHeight = 45.75
$ws.Rows.Item(17).RowHeight = 45.75
$ws.Rows.Item(18).RowHeight = 45.75

# --- Sheet view: scroll so row 19 is at the top, select D32 ---
$ws.Range("D32").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
